$d = $word.ActiveDocument

$replacements = @(
    @{old="486×5="; new="824×4="},
    @{old="837×7="; new="257×8="},
    @{old="486×9="; new="114×3="},
    @{old="750×4="; new="703×9="},
    @{old="981×8="; new="450×2="},
    @{old="478×3="; new="842×5="},
    @{old="609×7="; new="878×4="},
    @{old="702×2="; new="210×4="},
    @{old="577×5="; new="196×9="},
    @{old="719×8="; new="411×8="},
    @{old="612×3="; new="330×8="},
    @{old="137×9="; new="750×6="},
    @{old="141×9="; new="531×5="},
    @{old="371×2="; new="415×4="},
    @{old="308×9="; new="433×3="},
    @{old="714×8="; new="646×6="},
    @{old="971×6="; new="400×9="},
    @{old="256×5="; new="254×8="},
    @{old="640×5="; new="651×4="},
    @{old="206×4="; new="636×9="},
    @{old="416×2="; new="652×4="},
    @{old="720×5="; new="986×2="},
    @{old="399×4="; new="107×8="},
    @{old="830×6="; new="614×2="},
    @{old="730×3="; new="841×6="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
